# Redo the results for Random Forest and KNN
# (also updates the SVM row and re-titles the first block; clears the
# Decision Tree / Naive Bayes / AdaBoost rows on both sheets, since those
# results are now stale.)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Dep or Non-Dep Without LOSOCV"
$ws2 = $wb.Worksheets.Item(2)   # "Dep or Non-Dep With LOSOCV"

# ---------------------------------------------------------------------------
# Section titles (sheet 1 only)
# ---------------------------------------------------------------------------
$ws1.Range("A3").Value = "Depression/Non-Depression Classification with all temporal features"
$ws1.Range("A13").Value = "Depression/Non-Depression Classification with feature selection"

# ---------------------------------------------------------------------------
# Sheet 1, table 1 (rows 5-10: Random Forest, KNN, SVM, Decision Tree, Naive Bayes, AdaBoost)
# ---------------------------------------------------------------------------
$ws1.Cells.Item(5, 2).Value = 0.69144058081772997
$ws1.Cells.Item(5, 3).Value = 0.60032102728731895
$ws1.Cells.Item(5, 4).Value = 0.40107238605898099
$ws1.Cells.Item(5, 5).Value = 0.48087431693989002
$ws1.Cells.Item(5, 6).Value = 0.62662702116840396

$ws1.Cells.Item(6, 2).Value = 0.647688192586931
$ws1.Cells.Item(6, 3).Value = 0.50735809390329301
$ws1.Cells.Item(6, 4).Value = 0.38820375335120599
$ws1.Cells.Item(6, 5).Value = 0.43985419198055797
$ws1.Cells.Item(6, 6).Value = 0.58976824651828597

$ws1.Cells.Item(7, 2).Value = 0.64367596484524203
$ws1.Cells.Item(7, 3).Value = 0
$ws1.Cells.Item(7, 4).Value = 0
$ws1.Cells.Item(7, 5).Value = 0
$ws1.Cells.Item(7, 6).Value = 0.5
$ws1.Cells.Item(7, 6).NumberFormat = "0.0"

$ws1.Range("B8:F8").ClearContents()
$ws1.Range("B9:F9").ClearContents()
$ws1.Range("B10:F10").ClearContents()

# ---------------------------------------------------------------------------
# Sheet 1, table 2 (rows 15-20)
# ---------------------------------------------------------------------------
$ws1.Cells.Item(15, 2).Value = 0.67233473442873504
$ws1.Cells.Item(15, 3).Value = 0.55742725880551303
$ws1.Cells.Item(15, 4).Value = 0.39034852546916798
$ws1.Cells.Item(15, 5).Value = 0.45916114790286899
$ws1.Cells.Item(15, 6).Value = 0.60939213153838301

$ws1.Cells.Item(16, 2).Value = 0.647879251050821
$ws1.Cells.Item(16, 3).Value = 0.50784593437945702
$ws1.Cells.Item(16, 4).Value = 0.38176943699731902
$ws1.Cells.Item(16, 5).Value = 0.43587389041934499
$ws1.Cells.Item(16, 6).Value = 0.58848044423329804

$ws1.Range("B17:F17").ClearContents()
$ws1.Range("B18:F18").ClearContents()
$ws1.Range("B19:F19").ClearContents()
$ws1.Range("B20:F20").ClearContents()

# ---------------------------------------------------------------------------
# Sheet 2, table 1 (rows 5-10)
# ---------------------------------------------------------------------------
$ws2.Cells.Item(5, 2).Value = 0.60813489258062603
$ws2.Cells.Item(5, 3).Value = 0.41818181818181799
$ws2.Cells.Item(5, 4).Value = 0.13853297896702699
$ws2.Cells.Item(5, 5).Value = 0.20185682737924099
$ws2.Cells.Item(5, 6).Value = 0.60813489258062603

$ws2.Cells.Item(6, 2).Value = 0.54367002752932303
$ws2.Cells.Item(6, 3).Value = 0.41818181818181799
$ws2.Cells.Item(6, 4).Value = 0.14819598480382101
$ws2.Cells.Item(6, 5).Value = 0.21754541989269699
$ws2.Cells.Item(6, 6).Value = 0.54367002752932303
$ws2.Cells.Item(6, 6).NumberFormat = "0.00000"

$ws2.Range("B7:F7").ClearContents()
$ws2.Range("B8:F8").ClearContents()
$ws2.Range("B9:F9").ClearContents()
$ws2.Range("B10:F10").ClearContents()

# ---------------------------------------------------------------------------
# Sheet 2, table 2 (rows 15-20)
# ---------------------------------------------------------------------------
$ws2.Cells.Item(15, 2).Value = 0.595048989455223
$ws2.Cells.Item(15, 3).Value = 0.41818181818181799
$ws2.Cells.Item(15, 4).Value = 0.14251313627492901
$ws2.Cells.Item(15, 5).Value = 0.20779761676101999
$ws2.Cells.Item(15, 6).Value = 0.595048989455223

$ws2.Cells.Item(16, 2).Value = 0.54687729314767897
$ws2.Cells.Item(16, 3).Value = 0.41818181818181799
$ws2.Cells.Item(16, 4).Value = 0.149612988247835
$ws2.Cells.Item(16, 5).Value = 0.218975170924364
$ws2.Cells.Item(16, 6).Value = 0.54687729314767897

$ws2.Range("B17:F17").ClearContents()
$ws2.Range("B18:F18").ClearContents()
$ws2.Range("B19:F19").ClearContents()
$ws2.Range("B20:F20").ClearContents()
